$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.92361284117466
$ws.Range("C2").Value = 10.54777194214716
$ws.Range("D2").Value = 5.973310286292193
$ws.Range("E2").Value = 16.50105917307264
$ws.Range("G2").Value = 43.43456739963639
$ws.Range("H2").Value = 16.37313746056071
$ws.Range("B3").Value = 15.23726017800609
$ws.Range("C3").Value = 9.852321343010539
$ws.Range("D3").Value = 5.854069282555279
$ws.Range("E3").Value = 15.55777290136715
$ws.Range("G3").Value = 42.1053213618592
$ws.Range("H3").Value = 16.24519574023489
$ws.Range("B4").Value = 14.80495104501015
$ws.Range("C4").Value = 9.398936252868783
$ws.Range("D4").Value = 5.781693534853252
$ws.Range("E4").Value = 14.9546961589656
$ws.Range("G4").Value = 41.28548885698525
$ws.Range("H4").Value = 16.17115921909263
$ws.Range("B5").Value = 14.62634306243389
$ws.Range("C5").Value = 9.207523888391128
$ws.Range("D5").Value = 5.752454022539622
$ws.Range("E5").Value = 14.70319893162104
$ws.Range("G5").Value = 40.95103312704209
$ws.Range("H5").Value = 16.14214273136993
$ws.Range("B6").Value = 14.59654750520033
$ws.Range("C6").Value = 9.175337477102685
$ws.Range("D6").Value = 5.747615442356307
$ws.Range("E6").Value = 14.66110036315606
$ws.Range("G6").Value = 40.89549158169497
$ws.Range("H6").Value = 16.13739472140733
$ws.Range("B7").Value = 14.80255172727553
$ws.Range("C7").Value = 9.396381770561028
$ws.Range("D7").Value = 5.781298114094451
$ws.Range("E7").Value = 14.9513272185959
$ws.Range("G7").Value = 41.28097900330539
$ws.Range("H7").Value = 16.17076319902245
$ws.Range("B8").Value = 15.68938604657247
$ws.Range("C8").Value = 10.31344150268401
$ws.Range("D8").Value = 5.932047833800303
$ws.Range("E8").Value = 16.1809196962788
$ws.Range("G8").Value = 42.97732961778048
$ws.Range("H8").Value = 16.32809614817481
$ws.Range("B9").Value = 17.33073920920203
$ws.Range("C9").Value = 11.90352378336458
$ws.Range("D9").Value = 6.232437880512101
$ws.Range("E9").Value = 18.48617259078857
$ws.Range("G9").Value = 46.25187376538107
$ws.Range("H9").Value = 16.67161459347078
$ws.Range("B10").Value = 18.46416310473851
$ws.Range("C10").Value = 12.94620947552828
$ws.Range("D10").Value = 6.45357791833379
$ws.Range("E10").Value = 20.152354983554
$ws.Range("G10").Value = 48.59778943477246
$ws.Range("H10").Value = 16.94400391819327
$ws.Range("B11").Value = 18.96191511224383
$ws.Range("C11").Value = 13.39359950812954
$ws.Range("D11").Value = 6.553780859266603
$ws.Range("E11").Value = 20.86893411840885
$ws.Range("G11").Value = 49.64678310004593
$ws.Range("H11").Value = 17.07194267788611
$ws.Range("B12").Value = 19.14769067522225
$ws.Range("C12").Value = 13.55916573058657
$ws.Range("D12").Value = 6.591629420312642
$ws.Range("E12").Value = 21.13437707978427
$ws.Range("G12").Value = 50.0410129313067
$ws.Range("H12").Value = 17.12093930088309
$ws.Range("B13").Value = 19.10780329974009
$ws.Range("C13").Value = 13.52367898733625
$ws.Range("D13").Value = 6.583482934616441
$ws.Range("E13").Value = 21.07747140573046
$ws.Range("G13").Value = 49.95624773781547
$ws.Range("H13").Value = 17.11036307838117
$ws.Range("B14").Value = 18.97725409421079
$ws.Range("C14").Value = 13.40729778655098
$ws.Range("D14").Value = 6.556896827453254
$ws.Range("E14").Value = 20.89089048838735
$ws.Range("G14").Value = 49.67927886207434
$ws.Range("H14").Value = 17.07596283419317
$ws.Range("B15").Value = 18.89693176843526
$ws.Range("C15").Value = 13.3355101302511
$ws.Range("D15").Value = 6.540598419536077
$ws.Range("E15").Value = 20.77583581073895
$ws.Range("G15").Value = 49.50922571121736
$ws.Range("H15").Value = 17.05496226936613
$ws.Range("B16").Value = 18.43126170911249
$ws.Range("C16").Value = 12.91643142643522
$ws.Range("D16").Value = 6.44701806580672
$ws.Range("E16").Value = 20.10469566929626
$ws.Range("G16").Value = 48.5288354861048
$ws.Range("H16").Value = 16.93572131486428
$ws.Range("B17").Value = 18.14090583643773
$ws.Range("C17").Value = 12.65245786551281
$ws.Range("D17").Value = 6.389479793570627
$ws.Range("E17").Value = 19.68240396128009
$ws.Range("G17").Value = 47.92245474475822
$ws.Range("H17").Value = 16.86358135853085
$ws.Range("B18").Value = 17.97222803258557
$ws.Range("C18").Value = 12.49809086689045
$ws.Range("D18").Value = 6.356350372970267
$ws.Range("E18").Value = 19.43561610079631
$ws.Range("G18").Value = 47.57199058543674
$ws.Range("H18").Value = 16.82246934792222
$ws.Range("B19").Value = 17.91483414019774
$ws.Range("C19").Value = 12.44538854952199
$ws.Range("D19").Value = 6.345128554094193
$ws.Range("E19").Value = 19.35138747389338
$ws.Range("G19").Value = 47.45305162446961
$ws.Range("H19").Value = 16.80861585886944
$ws.Range("B20").Value = 18.1719889693767
$ws.Range("C20").Value = 12.68082064714636
$ws.Range("D20").Value = 6.395608747443685
$ws.Range("E20").Value = 19.72776065420631
$ws.Range("G20").Value = 47.98718284761224
$ws.Range("H20").Value = 16.87122156393976
$ws.Range("B21").Value = 19.01567422313103
$ws.Range("C21").Value = 13.44158609503085
$ws.Range("D21").Value = 6.564708723423128
$ws.Range("E21").Value = 20.94585390603109
$ws.Range("G21").Value = 49.76071569110167
$ws.Range("H21").Value = 17.08605236691134
$ws.Range("B22").Value = 19.55120943140147
$ws.Range("C22").Value = 13.91635596770753
$ws.Range("D22").Value = 6.674648341858315
$ws.Range("E22").Value = 21.7075186290274
$ws.Range("G22").Value = 50.90217474088803
$ws.Range("H22").Value = 17.2296418318634
$ws.Range("B23").Value = 19.26687808253874
$ws.Range("C23").Value = 13.66500773529066
$ws.Range("D23").Value = 6.616036744970632
$ws.Range("E23").Value = 21.30414102973097
$ws.Range("G23").Value = 50.29468962087594
$ws.Range("H23").Value = 17.15272427076185
$ws.Range("B24").Value = 18.15794171938526
$ws.Range("C24").Value = 12.66800594950329
$ws.Range("D24").Value = 6.392837999820713
$ws.Range("E24").Value = 19.70726737642092
$ws.Range("G24").Value = 47.95792499794161
$ws.Range("H24").Value = 16.86776629659516
$ws.Range("B25").Value = 16.89862309664548
$ws.Range("C25").Value = 11.49551508760758
$ws.Range("D25").Value = 6.150922776395749
$ws.Range("E25").Value = 17.83563121502668
$ws.Range("G25").Value = 45.37464151684937
$ws.Range("H25").Value = 16.57506646223355
